# "changed main to 4d cube"
# The search-results table went from 3 data rows (unidirectional s-t,
# unidirectional t-s, bidirectional) down to a single bidirectional-only
# result row, with new timing/expansion numbers and a new [g_F,g_B] value
# for the 4d-cube run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two now-obsolete result rows (old rows 3 and 4); this shrinks
# the used range down to A1:G2.
$ws.Rows("3:4").Delete()

# Update the single remaining result row (row 2) with the new 4d-cube values.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "bidirectional"
$ws.Range("C2").Value = 525
$ws.Range("D2").Value = 200
$ws.Range("E2").Value = 22
$ws.Range("F2").Value = "[7,0]"
$ws.Range("G2").Value = "file_path_here"
